$d = $word.ActiveDocument

# Locate the "Dewalt 20 Volt Power Tools" table: identify it by its last
# existing row, which is the Model Number "DCB094K" / "USB Charging Kit
# (with USB-C port and charger)" row that the new row is appended after.
$target = $null
for ($i = 1; $i -le $d.Tables.Count; $i++) {
    $candidate = $d.Tables.Item($i)
    $lastRow = $candidate.Rows.Count
    $modelNumber = $candidate.Cell($lastRow, 1).Range.Text.TrimEnd([char]7, [char]13)
    if ($modelNumber -eq "DCB094K") {
        $target = $candidate
        break
    }
}

# Append a new row for the DXSP190681 4 gal. Cordless Backpack Sprayer.
$newRow = $target.Rows.Add()
$newRow.Cells.Item(1).Range.Text = "DXSP190681"
$newRow.Cells.Item(2).Range.Text = "4 gal. Cordless Backpack Sprayer"
